# "Update countries & provincias Spain"
#
# The source list (sheet "Pais") is sorted by total cases. A COVID-19 data
# refresh (new timestamp "...a las 11:52") caused four countries - Kuwait,
# Albania, Fiyi and Botsuana - to move up in the ranking. Re-sorting
# shifted every row between each country's old and new rank down by one,
# so both the country name (column A) and the statistics (columns B-H:
# Casos totales, Nuevos casos, Casos activos, Recuperados, Casos
# criticos, Muertes hoy, Muertes) had to be rewritten for the affected
# rows. A handful of other rows only received refreshed statistics
# without changing rank/position.
#
# Below, every touched cell is set to its final value explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - footer timestamp sentence
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 11:52"

# Row 5 - Espana (rank unchanged, stats refreshed)
$ws.Range("B5").Value = 195944
$ws.Range("C5").Value = 1528
$ws.Range("D5").Value = 77357
$ws.Range("E5").Value = 97948

# Row 14 - Belgica (rank unchanged, stats refreshed)
$ws.Range("B14").Value = 38496
$ws.Range("C14").Value = 1313
$ws.Range("D14").Value = 8757
$ws.Range("E14").Value = 24056
$ws.Range("F14").Value = 1081
$ws.Range("G14").Value = 230
$ws.Range("H14").Value = 5683

# Row 60 - Tailandia (rank unchanged, stats refreshed)
$ws.Range("D60").Value = 457
$ws.Range("E60").Value = 1860
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 61

# Rows 63-66 - Kuwait moves up to rank 67, cascading Croacia/Barein/Islandia down one row
$ws.Range("A63").Value = "Kuwait"
$ws.Range("B63").Value = 1915
$ws.Range("C63").Value = 164
$ws.Range("D63").Value = 305
$ws.Range("E63").Value = 1603
$ws.Range("F63").Value = 38
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 7

$ws.Range("A64").Value = "Croacia"
$ws.Range("B64").Value = 1832
$ws.Range("D64").Value = 615
$ws.Range("E64").Value = 1178
$ws.Range("F64").Value = 27
$ws.Range("H64").Value = 39

$ws.Range("A65").Value = "Barein"
$ws.Range("B65").Value = 1773
$ws.Range("D65").Value = 755
$ws.Range("E65").Value = 1011
$ws.Range("H65").Value = 7

$ws.Range("A66").Value = "Islandia"
$ws.Range("B66").Value = 1760
$ws.Range("D66").Value = 1291
$ws.Range("E66").Value = 460
$ws.Range("F66").Value = 3
$ws.Range("H66").Value = 9

# Row 93 - Libano (rank unchanged, stats refreshed)
$ws.Range("B93").Value = 673
$ws.Range("C93").Value = 1
$ws.Range("E93").Value = 553
$ws.Range("F93").Value = 27

# Rows 97-98 - Albania moves up to rank 101, cascading Kirguistan down one row
$ws.Range("A97").Value = "Albania"
$ws.Range("B97").Value = 562
$ws.Range("C97").Value = 14
$ws.Range("D97").Value = 314
$ws.Range("E97").Value = 222
$ws.Range("H97").Value = 26

$ws.Range("A98").Value = "Kirguistan"
$ws.Range("B98").Value = 554
$ws.Range("C98").Value = 48
$ws.Range("D98").Value = 133
$ws.Range("E98").Value = 416
$ws.Range("H98").Value = 5

# Rows 179-188 - Botsuana moves up to rank 183, cascading Laos/Timor Oriental/
# Belice/Nueva Caledonia/Fiyi(unchanged at 184)/Islas Virgenes de los Estados
# Unidos/Malaui/Namibia/Dominica down one row each
$ws.Range("A179").Value = "Botsuana"
$ws.Range("B179").Value = 20
$ws.Range("C179").Value = 5
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 19
$ws.Range("H179").Value = 1

$ws.Range("A180").Value = "Laos"
$ws.Range("B180").Value = 19
$ws.Range("D180").Value = 2

$ws.Range("A181").Value = "Timor Oriental"
$ws.Range("D181").Value = 1
$ws.Range("E181").Value = 17
$ws.Range("F181").Value = 0
$ws.Range("H181").Value = 0

$ws.Range("A182").Value = "Belice"
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 16
$ws.Range("H182").Value = 2

$ws.Range("A183").Value = "Nueva Caledonia"
$ws.Range("B183").Value = 18
$ws.Range("D183").Value = 15
$ws.Range("E183").Value = 3
$ws.Range("F183").Value = 1

# Row 184 (Fiyi) is unchanged - it is the resync point of the cascade

$ws.Range("A185").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 17
$ws.Range("F185").Value = 0
$ws.Range("H185").Value = 0

$ws.Range("A186").Value = "Malaui"
$ws.Range("B186").Value = 17
$ws.Range("D186").Value = 3
$ws.Range("E186").Value = 12
$ws.Range("F186").Value = 1
$ws.Range("H186").Value = 2

$ws.Range("A187").Value = "Namibia"
$ws.Range("D187").Value = 6
$ws.Range("E187").Value = 10

$ws.Range("A188").Value = "Dominica"
$ws.Range("B188").Value = 16
$ws.Range("D188").Value = 8
$ws.Range("E188").Value = 8
$ws.Range("H188").Value = 0
